$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 8-10 (previously "ECs" as sending cluster); the
# new TPM data only has FAPs/MuSCs as sending clusters (6 data rows total).
$ws.Range("A8:A10").EntireRow.Delete() | Out-Null

# New per-row data (updated TPM values), columns A-D are text, E-T numeric.
$data = @(
  @{ Row=2;  A="FAPs";  B="Edn3"; C="Ednra"; D="ECs";   E=1; F=0.3333333333333333; G=0.06003666666666666; H=0.18011;           I=0.007162610180657564; J=0.007162610180657565; K=3; L=1; M=1.164012;           N=3.492036;   O=0.02222380689314669; P=0.02222380689314669; Q=0.06988340043999998; R=0.6289506039599999;   S=0.0001591804655058202; T=0.0001591804655058202 }
  @{ Row=3;  A="FAPs";  B="Edn3"; C="Ednra"; D="FAPs";  E=1; F=0.3333333333333333; G=0.06003666666666666; H=0.18011;           I=0.007162610180657564; J=0.007162610180657565; K=3; L=1; M=14.70158366666666;   N=44.10475099999999; O=0.2806888214480945;  P=0.2806888214480945; Q=0.8826340780677776;  R=7.943706702609998;    S=0.002010464610100895; T=0.002010464610100896 }
  @{ Row=4;  A="FAPs";  B="Edn3"; C="Ednra"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.06003666666666666; H=0.18011;           I=0.007162610180657564; J=0.007162610180657565; K=3; L=1; M=36.51120933333333;   N=109.533628;         O=0.6970873716587588;  P=0.6970873716587588; Q=2.192011304342222;   R=19.72810173908;       S=0.004992965105050848; T=0.004992965105050849 }
  @{ Row=5;  A="MuSCs"; B="Edn3"; C="Ednra"; D="ECs";   E=3; F=1;                   G=8.321916999999999;   H=24.965751;          I=0.9928373898193424;  J=0.9928373898193424;  K=3; L=1; M=1.164012;           N=3.492036;   O=0.02222380689314669; P=0.02222380689314669; Q=9.686811251003997;   R=87.18130125903598;    S=0.02206462642764087;  T=0.02206462642764087 }
  @{ Row=6;  A="MuSCs"; B="Edn3"; C="Ednra"; D="FAPs";  E=3; F=1;                   G=8.321916999999999;   H=24.965751;          I=0.9928373898193424;  J=0.9928373898193424;  K=3; L=1; M=14.70158366666666;   N=44.10475099999999; O=0.2806888214480945;  P=0.2806888214480945; Q=122.3453590425556;   R=1101.108231383001;    S=0.2786783568379936;   T=0.2786783568379936 }
  @{ Row=7;  A="MuSCs"; B="Edn3"; C="Ednra"; D="MuSCs"; E=3; F=1;                   G=8.321916999999999;   H=24.965751;          I=0.9928373898193424;  J=0.9928373898193424;  K=3; L=1; M=36.51120933333333;   N=109.533628;         O=0.6970873716587588;  P=0.6970873716587588; Q=303.8432536416253;   R=2734.589282774628;    S=0.6920944065537079;   T=0.6920944065537079 }
)

foreach ($r in $data) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value2  = $r.A
  $ws.Cells.Item($row, 2).Value2  = $r.B
  $ws.Cells.Item($row, 3).Value2  = $r.C
  $ws.Cells.Item($row, 4).Value2  = $r.D
  $ws.Cells.Item($row, 5).Value2  = $r.E
  $ws.Cells.Item($row, 6).Value2  = $r.F
  $ws.Cells.Item($row, 7).Value2  = $r.G
  $ws.Cells.Item($row, 8).Value2  = $r.H
  $ws.Cells.Item($row, 9).Value2  = $r.I
  $ws.Cells.Item($row, 10).Value2 = $r.J
  $ws.Cells.Item($row, 11).Value2 = $r.K
  $ws.Cells.Item($row, 12).Value2 = $r.L
  $ws.Cells.Item($row, 13).Value2 = $r.M
  $ws.Cells.Item($row, 14).Value2 = $r.N
  $ws.Cells.Item($row, 15).Value2 = $r.O
  $ws.Cells.Item($row, 16).Value2 = $r.P
  $ws.Cells.Item($row, 17).Value2 = $r.Q
  $ws.Cells.Item($row, 18).Value2 = $r.R
  $ws.Cells.Item($row, 19).Value2 = $r.S
  $ws.Cells.Item($row, 20).Value2 = $r.T
}
